# Insert a new record row before the existing row 133 in the "Maracuyá" daily
# price sheet (weekly refresh: one new observation added, all subsequent rows
# shift down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 133:164 down to 134:165, opening up a blank row 133.
$ws.Rows.Item(133).Insert()

# Populate the new row 133 with the new record. The descriptive columns
# (market/region/product/category/unit/origin) are identical to every other
# row in this block; only the date, quality grade, volume, min/avg/max price
# and $/kg columns differ per record.
$ws.Range("A133").Value = 1
$ws.Range("B133").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C133").Value = "Arica y Parinacota"
$ws.Range("D133").Value = 45005
$ws.Range("E133").Value = 15
$ws.Range("F133").Value = "Fruta"
$ws.Range("G133").Value = 100108
$ws.Range("H133").Value = "Tropicales y subtropicales"
$ws.Range("I133").Value = 100108003
$ws.Range("J133").Value = "Maracuyá"
$ws.Range("K133").Value = "Sin especificar"
$ws.Range("L133").Value = "Primera"
$ws.Range("M133").Value = 170
$ws.Range("N133").Value = 20000
$ws.Range("O133").Value = 22000
$ws.Range("P133").Value = 21000
$ws.Range("Q133").Value = "$/caja 20 kilos"
$ws.Range("R133").Value = "Región de Arica y Parinacota"
$ws.Range("S133").Value = 1050
$ws.Range("T133").Value = 20
